$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("C14").Value = 18527007950.849
$ws.Range("E14").Value = 36838371246.885
$ws.Range("G14").Value = 15289114716.375
$ws.Range("I14").Value = 58387627777.395
$ws.Range("J14").Value = 1.16312664126225
$ws.Range("K14").Value = 0.9209130678159025

# Row 17
$ws.Range("C17").Value = 10345869717.29
$ws.Range("E17").Value = 44190259843.321
$ws.Range("G17").Value = 28183349264.257
$ws.Range("I17").Value = 60197170422.38499
$ws.Range("J17").Value = 1.547178827538615
$ws.Range("K17").Value = 0.9202190473452521

# Row 18
$ws.Range("C18").Value = 9485654695.276
$ws.Range("E18").Value = 37542603376.893
$ws.Range("G18").Value = 23580962951.128
$ws.Range("I18").Value = 51504243802.658
$ws.Range("J18").Value = 1.471868929902973
$ws.Range("K18").Value = 0.9437932288052308

# Row 19
$ws.Range("C19").Value = 7578675569.81
$ws.Range("E19").Value = 24012465084.335
$ws.Range("G19").Value = 14674539375.694
$ws.Range("I19").Value = 33350390792.976
$ws.Range("J19").Value = 1.232131606984082
$ws.Range("K19").Value = 0.9217923022745748

# Row 20
$ws.Range("C20").Value = 9447401489.256
$ws.Range("E20").Value = 33624563124.17
$ws.Range("G20").Value = 22245635203.469
$ws.Range("I20").Value = 45003491044.871
$ws.Range("J20").Value = 1.204450550094819
$ws.Range("K20").Value = 1.025207727741531

# Row 21
$ws.Range("C21").Value = 14260494508.633
$ws.Range("E21").Value = 52905676359.56
$ws.Range("G21").Value = 35102741117.612
$ws.Range("I21").Value = 70708611601.508
$ws.Range("J21").Value = 1.248409389391825
$ws.Range("K21").Value = 0.9691767545540825

# Row 22
$ws.Range("C22").Value = 9152383958.514
$ws.Range("E22").Value = 25988479965.46
$ws.Range("G22").Value = 25988479965.46
$ws.Range("I22").Value = 25988479965.46

